$wb = $excel.ActiveWorkbook

# This script updates static market-data values (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns) that were
# refreshed by the scheduled data-update runner. Each sheet corresponds to a
# crafting job (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 5208593.5
$ws.Range("I12").Value = 6944625.5
$ws.Range("K12").Value = 6944625.5
$ws.Range("M12").Value = -6944455.5
$ws.Range("H38").Value = 3031.7144
$ws.Range("I38").Value = 156.625
$ws.Range("K38").Value = 469.875
$ws.Range("M38").Value = -97.875
$ws.Range("H64").Value = 9803.799999999999
$ws.Range("I64").Value = 7234.2
$ws.Range("K64").Value = 7234.2
$ws.Range("M64").Value = -6986.2
$ws.Range("H67").Value = 9803.799999999999
$ws.Range("I67").Value = 7234.2
$ws.Range("K67").Value = 7234.2
$ws.Range("M67").Value = -6376.2
$ws.Range("H80").Value = 326.25
$ws.Range("I80").Value = 248.21428
$ws.Range("J80").Value = 435.5
$ws.Range("K80").Value = 744.64284
$ws.Range("L80").Value = 1306.5
$ws.Range("M80").Value = 253.35716
$ws.Range("N80").Value = -3302.5
$ws.Range("H83").Value = 326.25
$ws.Range("I83").Value = 248.21428
$ws.Range("J83").Value = 435.5
$ws.Range("K83").Value = 2233.92852
$ws.Range("L83").Value = 3919.5
$ws.Range("M83").Value = 2758.07148
$ws.Range("N83").Value = -13903.5
$ws.Range("H111").Value = 2483.2727
$ws.Range("I111").Value = 2533.375
$ws.Range("J111").Value = 2349.6667
$ws.Range("K111").Value = 7600.125
$ws.Range("L111").Value = 7049.000100000001
$ws.Range("M111").Value = -4533.125
$ws.Range("N111").Value = -13183.0001
$ws.Range("H113").Value = 6491.1113
$ws.Range("J113").Value = 7686.5
$ws.Range("L113").Value = 7686.5
$ws.Range("N113").Value = -14194.5
$ws.Range("H132").Value = 1891.58
$ws.Range("I132").Value = 1812.1522
$ws.Range("K132").Value = 5436.4566
$ws.Range("M132").Value = -2906.4566
$ws.Range("H133").Value = 141651.67
$ws.Range("J133").Value = 141651.67
$ws.Range("L133").Value = 141651.67
$ws.Range("N133").Value = -151771.67
$ws.Range("H137").Value = 2223724.5
$ws.Range("I137").Value = 1022.5294
$ws.Range("K137").Value = 3067.5882
$ws.Range("M137").Value = -517.5882000000001
$ws.Range("H138").Value = 6071.603
$ws.Range("J138").Value = 4998.795
$ws.Range("L138").Value = 14996.385
$ws.Range("N138").Value = -25276.385

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 198955.4
$ws.Range("I32").Value = 202900.62
$ws.Range("K32").Value = 202900.62
$ws.Range("M32").Value = -202613.62
$ws.Range("H45").Value = 81948.08
$ws.Range("I45").Value = 88202.25
$ws.Range("K45").Value = 88202.25
$ws.Range("M45").Value = -87825.25
$ws.Range("H74").Value = 892820.5
$ws.Range("I74").Value = 2545.611
$ws.Range("J74").Value = 3182098.8
$ws.Range("K74").Value = 2545.611
$ws.Range("L74").Value = 3182098.8
$ws.Range("M74").Value = -1671.611
$ws.Range("N74").Value = -3183846.8
$ws.Range("H77").Value = 892820.5
$ws.Range("I77").Value = 2545.611
$ws.Range("J77").Value = 3182098.8
$ws.Range("K77").Value = 12728.055
$ws.Range("L77").Value = 15910494
$ws.Range("M77").Value = -8360.055
$ws.Range("N77").Value = -15919230
$ws.Range("H110").Value = 747.6875
$ws.Range("I110").Value = 597.7
$ws.Range("K110").Value = 597.7
$ws.Range("M110").Value = 1447.3
$ws.Range("H122").Value = 2261.2
$ws.Range("I122").Value = 1584.3334
$ws.Range("K122").Value = 4753.0002
$ws.Range("M122").Value = -2303.0002
$ws.Range("H132").Value = 4949.1816
$ws.Range("I132").Value = 3573.5833
$ws.Range("K132").Value = 10720.7499
$ws.Range("M132").Value = -8190.749899999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1592.4736
$ws.Range("I20").Value = 1455.25
$ws.Range("K20").Value = 1455.25
$ws.Range("M20").Value = -1208.25
$ws.Range("H58").Value = 80000
$ws.Range("J58").Value = 80000
$ws.Range("L58").Value = 80000
$ws.Range("N58").Value = -80588
$ws.Range("H86").Value = 1518.7241
$ws.Range("I86").Value = 1374.8334
$ws.Range("J86").Value = 1754.1818
$ws.Range("K86").Value = 1374.8334
$ws.Range("L86").Value = 1754.1818
$ws.Range("M86").Value = -251.8334
$ws.Range("N86").Value = -4000.1818
$ws.Range("H89").Value = 1518.7241
$ws.Range("I89").Value = 1374.8334
$ws.Range("J89").Value = 1754.1818
$ws.Range("K89").Value = 6874.166999999999
$ws.Range("L89").Value = 8770.909
$ws.Range("M89").Value = -1258.166999999999
$ws.Range("N89").Value = -20002.909
$ws.Range("H94").Value = 1426.341
$ws.Range("I94").Value = 1436.3438
$ws.Range("K94").Value = 1436.3438
$ws.Range("M94").Value = -985.3438000000001
$ws.Range("H105").Value = 5036.6553
$ws.Range("I105").Value = 6469.05
$ws.Range("J105").Value = 1853.5555
$ws.Range("K105").Value = 6469.05
$ws.Range("L105").Value = 1853.5555
$ws.Range("M105").Value = -4722.05
$ws.Range("N105").Value = -5347.5555

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2769.2
$ws.Range("I31").Value = 1276.849
$ws.Range("J31").Value = 5240.9062
$ws.Range("K31").Value = 1276.849
$ws.Range("L31").Value = 5240.9062
$ws.Range("M31").Value = -981.8489999999999
$ws.Range("N31").Value = -5830.9062
$ws.Range("H34").Value = 2769.2
$ws.Range("I34").Value = 1276.849
$ws.Range("J34").Value = 5240.9062
$ws.Range("K34").Value = 1276.849
$ws.Range("L34").Value = 5240.9062
$ws.Range("M34").Value = -1074.849
$ws.Range("N34").Value = -5644.9062
$ws.Range("H134").Value = 3328.5454
$ws.Range("I134").Value = 2377.4
$ws.Range("J134").Value = 4121.1665
$ws.Range("K134").Value = 7132.200000000001
$ws.Range("L134").Value = 12363.4995
$ws.Range("M134").Value = -4597.200000000001
$ws.Range("N134").Value = -17433.4995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3086.2307
$ws.Range("I68").Value = 2741
$ws.Range("K68").Value = 8223
$ws.Range("M68").Value = -7412
$ws.Range("H71").Value = 3086.2307
$ws.Range("I71").Value = 2741
$ws.Range("K71").Value = 24669
$ws.Range("M71").Value = -20613
$ws.Range("H107").Value = 250000140
$ws.Range("I107").Value = 196
$ws.Range("K107").Value = 588
$ws.Range("M107").Value = 1332
$ws.Range("H109").Value = 5010.6665
$ws.Range("I109").Value = 1667
$ws.Range("K109").Value = 5001
$ws.Range("M109").Value = -3961

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 41561.688
$ws.Range("I49").Value = 29999.5
$ws.Range("J49").Value = 43213.43
$ws.Range("K49").Value = 29999.5
$ws.Range("L49").Value = 43213.43
$ws.Range("M49").Value = -29815.5
$ws.Range("N49").Value = -43581.43
$ws.Range("H80").Value = 116406.4
$ws.Range("I80").Value = 121350.88
$ws.Range("K80").Value = 121350.88
$ws.Range("M80").Value = -120352.88
$ws.Range("H83").Value = 116406.4
$ws.Range("I83").Value = 121350.88
$ws.Range("K83").Value = 606754.4
$ws.Range("M83").Value = -601762.4
$ws.Range("H126").Value = 1905.5625
$ws.Range("I126").Value = 1773.8462
$ws.Range("K126").Value = 5321.5386
$ws.Range("M126").Value = -2851.5386
$ws.Range("H132").Value = 11987478
$ws.Range("I132").Value = 3077.7856
$ws.Range("J132").Value = 30629878
$ws.Range("K132").Value = 9233.356800000001
$ws.Range("L132").Value = 91889634
$ws.Range("M132").Value = -6703.356800000001
$ws.Range("N132").Value = -91894694

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3938.0527
$ws.Range("I40").Value = 3526.6667
$ws.Range("K40").Value = 3526.6667
$ws.Range("M40").Value = -3390.6667
$ws.Range("H80").Value = 77479
$ws.Range("J80").Value = 77479
$ws.Range("L80").Value = 77479
$ws.Range("N80").Value = -79725
$ws.Range("H83").Value = 77479
$ws.Range("J83").Value = 77479
$ws.Range("L83").Value = 232437
$ws.Range("N83").Value = -243669
$ws.Range("H122").Value = 4749.38
$ws.Range("I122").Value = 3893.9656
$ws.Range("J122").Value = 5930.6665
$ws.Range("K122").Value = 11681.8968
$ws.Range("L122").Value = 17791.9995
$ws.Range("M122").Value = -9231.8968
$ws.Range("N122").Value = -22691.9995
$ws.Range("H132").Value = 7709.9
$ws.Range("I132").Value = 2406.3635
$ws.Range("J132").Value = 14192
$ws.Range("K132").Value = 7219.0905
$ws.Range("L132").Value = 42576
$ws.Range("M132").Value = -4689.0905
$ws.Range("N132").Value = -47636

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1325.6957
$ws.Range("I81").Value = 1057.7693
$ws.Range("J81").Value = 1674
$ws.Range("K81").Value = 2115.5386
$ws.Range("L81").Value = 3348
$ws.Range("M81").Value = -1054.5386
$ws.Range("N81").Value = -5470
$ws.Range("H84").Value = 1325.6957
$ws.Range("I84").Value = 1057.7693
$ws.Range("J84").Value = 1674
$ws.Range("K84").Value = 10577.693
$ws.Range("L84").Value = 16740
$ws.Range("M84").Value = -5273.692999999999
$ws.Range("N84").Value = -27348
$ws.Range("H122").Value = 4053.0417
$ws.Range("I122").Value = 4652.778
$ws.Range("K122").Value = 13958.334
$ws.Range("M122").Value = -11508.334
$ws.Range("H132").Value = 37979.17
$ws.Range("I132").Value = 55771.473
$ws.Range("J132").Value = 4173.8
$ws.Range("K132").Value = 167314.419
$ws.Range("L132").Value = 12521.4
$ws.Range("M132").Value = -164784.419
$ws.Range("N132").Value = -17581.4
